$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A35").Value = "Howl"
$ws.Range("B35").Value = "Curse"
$ws.Range("D35").Value = "Beast"
$ws.Range("F35").Value = "3 turns"
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = "SPR (Endurance)"
$ws.Range("I35").Value = 10
$ws.Range("J35").Value = "{\bf Werewolf Species spell. This spell can only be learned by werewolves} \\ Release an earsplitting, supernatural roar which causes all beings within 100m to perform a SPR Resist. Failure causes them to gain the {\it Terrified} status. "

$ws.Range("A36").Value = "Confound"
$ws.Range("B36").Value = "Curse"
$ws.Range("C36").Value = "lombus"
$ws.Range("D36").Value = "Instant"
$ws.Range("E36").Value = "Blue bolt"
$ws.Range("F36").Value = "2 turns"
$ws.Range("G36").Value = 1
$ws.Range("H36").Value = "POW"
$ws.Range("I36").Value = "CC"
$ws.Range("J36").Value = "The target suffers a 1-point penalty to all checks for the duration of the spell. "

$ws.Rows(35).RowHeight = 19.4

$ws.Range("A37").Select() | Out-Null
